# Updated cryptos list (prices / 1h volume %) to match the latest scrape.
# Note: several "Price" cells look numeric (e.g. 247.65) but must stay TEXT,
# exactly like the rest of column D already is. A leading apostrophe forces
# Excel to store them as text instead of silently coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.212.27'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '2.247.48'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''247.65'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').Value = '''76.13'
$ws.Range('E7').Value = '  +5.01%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '''0.631'
$ws.Range('E9').Value = '  -2.16%  '
$ws.Range('D10').Value = '''40.30'
$ws.Range('E10').Value = '  +3.40%  '
$ws.Range('D11').Value = '''0.0947'
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('E13').Value = '  -1.54%  '
$ws.Range('D14').Value = '2.584.64'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '''14.90'
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('D16').Value = '''0.861'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '2.256.62'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '42.213.46'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').Value = '0.0₃0976'
$ws.Range('E19').Value = '  -2.29%  '
$ws.Range('E20').Value = '  -2.39%  '
$ws.Range('D21').Value = '''71.62'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = '''2.23'
$ws.Range('E22').Value = '  -4.25%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '''231.96'
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E26').Value = '  -4.45%  '
$ws.Range('E27').Value = '  -4.56%  '
$ws.Range('D28').Value = '''7.13'
$ws.Range('E28').Value = '  +10.13%  '
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').Value = '''168.60'
$ws.Range('E30').Value = '  +0.74%  '
$ws.Range('D31').Value = '''20.56'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('D32').Value = '''0.0847'
$ws.Range('E32').Value = '  +5.33%  '
$ws.Range('D33').Value = '''32.78'
$ws.Range('E33').Value = '  +3.05%  '
$ws.Range('D34').Value = '''0.120'
$ws.Range('E34').Value = '  -6.62%  '
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').Value = '''4.53'
$ws.Range('E36').Value = '  -5.08%  '
$ws.Range('D37').Value = '''4.80'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').Value = '''0.0297'
$ws.Range('E38').Value = '  -3.65%  '
$ws.Range('D39').Value = '''13.24'
$ws.Range('E39').Value = '  -6.67%  '
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('E41').Value = '  -5.69%  '
$ws.Range('D42').Value = '''118.09'
$ws.Range('E42').Value = '  +21.42%  '
$ws.Range('E43').Value = '  -4.82%  '
$ws.Range('D44').Value = '''60.17'
$ws.Range('E44').Value = '  -3.15%  '
$ws.Range('D45').Value = '''8.72'
$ws.Range('E45').Value = '  -6.28%  '
$ws.Range('E46').Value = '  -2.74%  '
$ws.Range('D47').Value = '''0.997'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  -4.24%  '
$ws.Range('E49').Value = '  -1.83%  '
$ws.Range('D50').Value = '''4.29'
$ws.Range('E50').Value = '  -12.44%  '
$ws.Range('D51').Value = '''4.13'
$ws.Range('E51').Value = '  -2.70%  '
